$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3485.2646
$ws.Range("I64").Value = 3212.375
$ws.Range("J64").Value = 3569.2307
$ws.Range("K64").Value = 3212.375
$ws.Range("L64").Value = 3569.2307
$ws.Range("M64").Value = -2964.375
$ws.Range("N64").Value = -4065.2307
$ws.Range("H67").Value = 3485.2646
$ws.Range("I67").Value = 3212.375
$ws.Range("J67").Value = 3569.2307
$ws.Range("K67").Value = 3212.375
$ws.Range("L67").Value = 3569.2307
$ws.Range("M67").Value = -2354.375
$ws.Range("N67").Value = -5285.2307
$ws.Range("H106").Value = 2738.3333
$ws.Range("I106").Value = 2886
$ws.Range("K106").Value = 2886
$ws.Range("M106").Value = -2255
$ws.Range("H129").Value = 1261.6364
$ws.Range("J129").Value = 1455.5
$ws.Range("L129").Value = 4366.5
$ws.Range("N129").Value = -14366.5
$ws.Range("H137").Value = 1333.4445
$ws.Range("I137").Value = 1105.2632
$ws.Range("J137").Value = 1875.375
$ws.Range("K137").Value = 3315.7896
$ws.Range("L137").Value = 5626.125
$ws.Range("M137").Value = -765.7896000000001
$ws.Range("N137").Value = -10726.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8469.364
$ws.Range("I32").Value = 8961.071
$ws.Range("J32").Value = 6174.7334
$ws.Range("K32").Value = 8961.071
$ws.Range("L32").Value = 6174.7334
$ws.Range("M32").Value = -8674.071
$ws.Range("N32").Value = -6748.7334
$ws.Range("H45").Value = 1234.6666
$ws.Range("I45").Value = 1085.3334
$ws.Range("J45").Value = 1533.3334
$ws.Range("K45").Value = 1085.3334
$ws.Range("L45").Value = 1533.3334
$ws.Range("M45").Value = -708.3334
$ws.Range("N45").Value = -2287.3334
$ws.Range("H61").Value = 2191.7144
$ws.Range("I61").Value = 1525.2354
$ws.Range("J61").Value = 5024.25
$ws.Range("K61").Value = 1525.2354
$ws.Range("L61").Value = 5024.25
$ws.Range("M61").Value = -1313.2354
$ws.Range("N61").Value = -5448.25
$ws.Range("H74").Value = 775.5833
$ws.Range("I74").Value = 798.0476
$ws.Range("J74").Value = 744.13336
$ws.Range("K74").Value = 798.0476
$ws.Range("L74").Value = 744.13336
$ws.Range("M74").Value = 75.95240000000001
$ws.Range("N74").Value = -2492.13336
$ws.Range("H77").Value = 775.5833
$ws.Range("I77").Value = 798.0476
$ws.Range("J77").Value = 744.13336
$ws.Range("K77").Value = 3990.238
$ws.Range("L77").Value = 3720.6668
$ws.Range("M77").Value = 377.7620000000002
$ws.Range("N77").Value = -12456.6668
$ws.Range("H102").Value = 252450
$ws.Range("J102").Value = 502500
$ws.Range("L102").Value = 502500
$ws.Range("N102").Value = -505744
$ws.Range("H122").Value = 2649.2307
$ws.Range("I122").Value = 3402.889
$ws.Range("K122").Value = 10208.667
$ws.Range("M122").Value = -7758.667000000001
$ws.Range("H132").Value = 5522.9756
$ws.Range("I132").Value = 7587.619
$ws.Range("J132").Value = 3355.1
$ws.Range("K132").Value = 22762.857
$ws.Range("L132").Value = 10065.3
$ws.Range("M132").Value = -20232.857
$ws.Range("N132").Value = -15125.3
$ws.Range("H136").Value = 2191.7144
$ws.Range("I136").Value = 1525.2354
$ws.Range("J136").Value = 5024.25
$ws.Range("K136").Value = 4575.706200000001
$ws.Range("L136").Value = 15072.75
$ws.Range("M136").Value = -2025.706200000001
$ws.Range("N136").Value = -20172.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9408.448
$ws.Range("I134").Value = 2871
$ws.Range("J134").Value = 18669.834
$ws.Range("K134").Value = 8613
$ws.Range("L134").Value = 56009.50199999999
$ws.Range("M134").Value = -6078
$ws.Range("N134").Value = -61079.50199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2297.7222
$ws.Range("I31").Value = 1902.6774
$ws.Range("J31").Value = 4747
$ws.Range("K31").Value = 1902.6774
$ws.Range("L31").Value = 4747
$ws.Range("M31").Value = -1607.6774
$ws.Range("N31").Value = -5337
$ws.Range("H34").Value = 2297.7222
$ws.Range("I34").Value = 1902.6774
$ws.Range("J34").Value = 4747
$ws.Range("K34").Value = 1902.6774
$ws.Range("L34").Value = 4747
$ws.Range("M34").Value = -1700.6774
$ws.Range("N34").Value = -5151
$ws.Range("H58").Value = 862554
$ws.Range("I58").Value = 1123463.5
$ws.Range("J58").Value = 1552.8
$ws.Range("K58").Value = 1123463.5
$ws.Range("L58").Value = 1552.8
$ws.Range("M58").Value = -1123260.5
$ws.Range("N58").Value = -1958.8
$ws.Range("H132").Value = 288824.03
$ws.Range("I132").Value = 338754.78
$ws.Range("J132").Value = 3505.5715
$ws.Range("K132").Value = 1016264.34
$ws.Range("L132").Value = 10516.7145
$ws.Range("M132").Value = -1013734.34
$ws.Range("N132").Value = -15576.7145
$ws.Range("H134").Value = 1233.1406
$ws.Range("I134").Value = 911.76
$ws.Range("J134").Value = 2380.9285
$ws.Range("K134").Value = 2735.28
$ws.Range("L134").Value = 7142.7855
$ws.Range("M134").Value = -200.2799999999997
$ws.Range("N134").Value = -12212.7855
$ws.Range("H136").Value = 862554
$ws.Range("I136").Value = 1123463.5
$ws.Range("J136").Value = 1552.8
$ws.Range("K136").Value = 3370390.5
$ws.Range("L136").Value = 4658.4
$ws.Range("M136").Value = -3367840.5
$ws.Range("N136").Value = -9758.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1187.7778
$ws.Range("I121").Value = 856
$ws.Range("J121").Value = 1453.2
$ws.Range("K121").Value = 2568
$ws.Range("L121").Value = 4359.6
$ws.Range("M121").Value = -1258
$ws.Range("N121").Value = -6979.6
$ws.Range("H129").Value = 2942627.2
$ws.Range("J129").Value = 4168504.8
$ws.Range("L129").Value = 12505514.4
$ws.Range("N129").Value = -12515514.4
$ws.Range("H133").Value = 4278.95
$ws.Range("I133").Value = 1801.909
$ws.Range("J133").Value = 7306.4443
$ws.Range("K133").Value = 5405.727000000001
$ws.Range("L133").Value = 21919.3329
$ws.Range("M133").Value = -345.7270000000008
$ws.Range("N133").Value = -32039.3329
$ws.Range("H134").Value = 4537.357
$ws.Range("I134").Value = 2796.3635
$ws.Range("J134").Value = 5663.8823
$ws.Range("K134").Value = 8389.0905
$ws.Range("L134").Value = 16991.6469
$ws.Range("M134").Value = -3319.0905
$ws.Range("N134").Value = -27131.6469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 28000
$ws.Range("J62").Value = 28000
$ws.Range("L62").Value = 28000
$ws.Range("N62").Value = -29372
$ws.Range("H65").Value = 28000
$ws.Range("J65").Value = 28000
$ws.Range("L65").Value = 84000
$ws.Range("N65").Value = -90864
$ws.Range("H102").Value = 2483.6365
$ws.Range("I102").Value = 2255.842
$ws.Range("J102").Value = 2792.7856
$ws.Range("K102").Value = 2255.842
$ws.Range("L102").Value = 2792.7856
$ws.Range("M102").Value = -633.8420000000001
$ws.Range("N102").Value = -6036.7856
$ws.Range("H132").Value = 2086.7222
$ws.Range("I132").Value = 1280.6154
$ws.Range("J132").Value = 4182.6
$ws.Range("K132").Value = 3841.8462
$ws.Range("L132").Value = 12547.8
$ws.Range("M132").Value = -1311.8462
$ws.Range("N132").Value = -17607.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1222.2222
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 1625
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 1625
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -2001
$ws.Range("H132").Value = 3387.3057
$ws.Range("I132").Value = 2976.2144
$ws.Range("J132").Value = 4826.125
$ws.Range("K132").Value = 8928.643199999999
$ws.Range("L132").Value = 14478.375
$ws.Range("M132").Value = -6398.643199999999
$ws.Range("N132").Value = -19538.375
$ws.Range("H136").Value = 2128.0186
$ws.Range("I136").Value = 1617.2142
$ws.Range("J136").Value = 3915.8333
$ws.Range("K136").Value = 4851.642599999999
$ws.Range("L136").Value = 11747.4999
$ws.Range("M136").Value = -2301.642599999999
$ws.Range("N136").Value = -16847.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4600
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -6048
$ws.Range("H65").Value = 4600
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -30240
$ws.Range("H132").Value = 920.1836499999999
$ws.Range("I132").Value = 683.88635
$ws.Range("J132").Value = 2999.6
$ws.Range("K132").Value = 2051.65905
$ws.Range("L132").Value = 8998.799999999999
$ws.Range("M132").Value = 478.3409499999998
$ws.Range("N132").Value = -14058.8
$ws.Range("H136").Value = 1653.3182
$ws.Range("I136").Value = 1713
$ws.Range("J136").Value = 400
$ws.Range("K136").Value = 5139
$ws.Range("L136").Value = 1200
$ws.Range("M136").Value = -2589
$ws.Range("N136").Value = -6300
